$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (shifts old I..O to J..P)
$ws.Columns("I:I").Insert()

# Match the new column's width to the neighboring column H (closest achievable)
$ws.Columns("I:I").ColumnWidth = 19.3

# New header for the inserted column
$ws.Range("I1").Value = "Porcentaje_avance_actividades"

# Fill in the new column with the "Porcentaje_avance_actividades" values
$ws.Range("I2").Value = 100
$ws.Range("I3").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("I6").Value = 90
$ws.Range("I7").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("I9").Value = 100
$ws.Range("I10").Value = 100
$ws.Range("I11").Value = 95
$ws.Range("I12").Value = 100
$ws.Range("I13").Value = 100
$ws.Range("I14").Value = 100
$ws.Range("I15").Value = 100
$ws.Range("I16").Value = 100
$ws.Range("I17").Value = 100
$ws.Range("I18").Value = 97
$ws.Range("I19").Value = 100
$ws.Range("I20").Value = 100
$ws.Range("I21").Value = 100
$ws.Range("I22").Value = 70
$ws.Range("I23").Value = 65
$ws.Range("I24").Value = 80
$ws.Range("I25").Value = 60
$ws.Range("I26").Value = 68
$ws.Range("I27").Value = 79
$ws.Range("I28").Value = 90
$ws.Range("I29").Value = 87
$ws.Range("I30").Value = 58
$ws.Range("I31").Value = 75

# Move the active selection to I1 (matches the saved workbook view state)
$ws.Range("I1").Select() | Out-Null
